# Updated input data for example run
#
# The_Vale.xlsx input-data refresh:
#  - labels_ExpROW: collapse the 6 category/"Export" label pairs down to a
#    single "Export" label (the ExpROW data itself loses its per-category
#    breakdown, see below), so the labels sheet shrinks to one row/column.
#  - ExpROW: columns B:F (per-category export-to-ROW breakdown) are folded
#    into column A as a single total per commodity row.
#  - VA: the "Import" rows for several small regions are consolidated into
#    a single aggregate row, shrinking the sheet by 5 rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) labels_ExpROW (sheet "labels_ExpROW"): keep only the "Export" label
# ---------------------------------------------------------------------
$wsLabels = $wb.Worksheets.Item("labels_ExpROW")

# A1 already holds a category name ("Elms"); B1 holds "Export" - keep
# only the "Export" label, then drop the now-unused rows/column.
$wsLabels.Range("A1").Value = $wsLabels.Range("B1").Value()
$wsLabels.Range("A2:A6").EntireRow.Delete()
$wsLabels.Columns.Item(2).Delete()

$wsLabels.Range("D15").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) ExpROW (sheet "ExpROW"): collapse B:F into a single total in A
# ---------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("ExpROW")

# Compute each row's total (A:F) via a helper formula in an unused
# column, then paste the computed totals back into column A as values.
for ($r = 1; $r -le 24; $r++) {
    $wsExp.Cells.Item($r, 8).Formula = "=SUM(A" + $r + ":F" + $r + ")"
}
$wb.Application.Calculate()
for ($r = 1; $r -le 24; $r++) {
    $total = $wsExp.Cells.Item($r, 8).Value()
    $wsExp.Cells.Item($r, 1).Value = $total
}
$wsExp.Range("H1:H24").Clear()

# Drop the now-redundant per-category columns B:F.
$wsExp.Range("B1:F24").Delete()

$wsExp.Range("A24").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) VA (sheet "VA"): merge rows 26,28,30,32,34,36 into one total row
# ---------------------------------------------------------------------
$wsVA = $wb.Worksheets.Item("VA")

# Total the six rows to be consolidated before any shifting happens.
$wsVA.Range("D1").Formula = "=A26+A28+A30+A32+A34+A36"
$wb.Application.Calculate()
$mergedTotal = $wsVA.Range("D1").Value()
$wsVA.Range("D1").Clear()

# Remove the six rows bottom-to-top so earlier row numbers stay valid.
$wsVA.Rows.Item(36).Delete()
$wsVA.Rows.Item(34).Delete()
$wsVA.Rows.Item(32).Delete()
$wsVA.Rows.Item(30).Delete()
$wsVA.Rows.Item(28).Delete()
$wsVA.Rows.Item(26).Delete()

# The remaining rows (old 27,29,31,33,35) are now 26..30; append the
# consolidated total as the new last row (31).
$wsVA.Range("A31").Value = $mergedTotal
$wsVA.Range("B31").Value = 0

# Make VA the active sheet/selection, matching the saved workbook state.
$wsVA.Activate() | Out-Null
$wsVA.Range("B29").Select() | Out-Null
